# Update the "Förändrad" date column (C2:C82) from 2023-09-03 (45172)
# to 2023-09-06 (45175) for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C82")
$range.Value = 45175
